$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-31 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-01 Sunday", 2)
$d.Content.Find.Execute("27÷3=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2)
$d.Content.Find.Execute("93÷5=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "83÷3=27, 2", 2)
$d.Content.Find.Execute("72÷8=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 2)
$d.Content.Find.Execute("17÷6=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "65÷5=13, 0", 2)
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷6=7, 4", 2)
$d.Content.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=12, 0", 2)
$d.Content.Find.Execute("71÷5=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=6, 0", 2)
$d.Content.Find.Execute("69÷3=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷8=7, 2", 2)
$d.Content.Find.Execute("39÷8=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "85÷5=17, 0", 2)
$d.Content.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "68÷7=9, 5", 2)
$d.Content.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2)
$d.Content.Find.Execute("69÷6=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=5, 4", 2)
$d.Content.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=13, 1", 2)
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=10, 5", 2)
$d.Content.Find.Execute("77÷2=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=12, 3", 2)
$d.Content.Find.Execute("65÷4=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=19, 1", 2)
$d.Content.Find.Execute("30÷4=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2)
$d.Content.Find.Execute("65÷7=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=6, 0", 2)
$d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷4=23, 3", 2)
$d.Content.Find.Execute("65÷3=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "43÷6=7, 1", 2)
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "24÷6=4, 0", 2)
$d.Content.Find.Execute("83÷4=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=15, 3", 2)
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=5, 1", 2)
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2)
$d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2)

Write-Host "Replacements complete"
